$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Row, $Col, $Val)
    $cell = $ws.Cells.Item($Row, $Col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $Val
    $cell.Style = $origStyle
}

Set-TextValue 2 4 '63.410.00'
Set-TextValue 2 5 '  +3.54%  '
Set-TextValue 3 4 '3.066.30'
Set-TextValue 3 5 '  +2.04%  '
Set-TextValue 4 5 '  +0.16%  '
Set-TextValue 5 4 '549.47'
Set-TextValue 5 5 '  +2.37%  '
Set-TextValue 6 4 '140.50'
Set-TextValue 6 5 '  +3.79%  '
Set-TextValue 7 5 '  -0.10%  '
Set-TextValue 8 4 '3.061.82'
Set-TextValue 8 5 '  +2.08%  '
Set-TextValue 9 4 '0.502'
Set-TextValue 9 5 '  +0.99%  '
Set-TextValue 10 4 '6.52'
Set-TextValue 10 5 '  +6.50%  '
Set-TextValue 11 5 '  +1.24%  '
Set-TextValue 12 5 '  +1.58%  '
Set-TextValue 13 5 '  +2.60%  '
Set-TextValue 14 4 '34.85'
Set-TextValue 14 5 '  +2.01%  '
Set-TextValue 15 4 '3.561.54'
Set-TextValue 15 5 '  +2.04%  '
Set-TextValue 16 4 '63.361.05'
Set-TextValue 16 5 '  +3.47%  '
Set-TextValue 17 4 '3.070.94'
Set-TextValue 17 5 '  +2.34%  '
Set-TextValue 18 5 '  -1.27%  '
Set-TextValue 19 5 '  +1.78%  '
Set-TextValue 20 4 '482.22'
Set-TextValue 20 5 '  +3.06%  '
Set-TextValue 21 5 '  +2.93%  '
Set-TextValue 22 5 '  -0.71%  '
Set-TextValue 23 4 '7.27'
Set-TextValue 23 5 '  +4.17%  '
Set-TextValue 24 4 '80.76'
Set-TextValue 24 5 '  +1.12%  '
Set-TextValue 25 4 '12.62'
Set-TextValue 25 5 '  +4.62%  '
Set-TextValue 26 4 '1.00'
Set-TextValue 26 5 '  -0.12%  '
Set-TextValue 27 5 '  +2.48%  '
Set-TextValue 28 5 '  -0.07%  '
Set-TextValue 29 5 '  +4.81%  '
Set-TextValue 30 5 '  +0.23%  '
Set-TextValue 31 4 '26.11'
Set-TextValue 31 5 '  +1.92%  '
Set-TextValue 32 5 '  -0.24%  '
Set-TextValue 33 4 '2.46'
Set-TextValue 33 5 '  +7.27%  '
Set-TextValue 34 4 '5.72'
Set-TextValue 34 5 '  +3.85%  '
Set-TextValue 35 4 '55.57'
Set-TextValue 35 5 '  -0.24%  '
Set-TextValue 36 5 '  +1.38%  '
Set-TextValue 37 4 '469.59'
Set-TextValue 37 5 '  +2.74%  '
Set-TextValue 38 4 '0.0820'
Set-TextValue 38 5 '  +3.95%  '
Set-TextValue 39 4 '0.0396'
Set-TextValue 39 5 '  +2.95%  '
Set-TextValue 40 4 '3.083.76'
Set-TextValue 40 5 '  -3.64%  '
Set-TextValue 41 5 '  -0.02%  '
Set-TextValue 42 5 '  +1.17%  '
Set-TextValue 43 5 '  +4.03%  '
Set-TextValue 44 4 '28.02'
Set-TextValue 44 5 '  +0.96%  '
Set-TextValue 45 5 '  +3.31%  '
Set-TextValue 46 5 '  -0.11%  '
Set-TextValue 47 5 '  +2.44%  '
Set-TextValue 48 5 '  +1.12%  '
Set-TextValue 49 4 '116.44'
Set-TextValue 49 5 '  -3.46%  '
Set-TextValue 50 4 '0.0₃0509'
Set-TextValue 50 5 '  +2.51%  '
Set-TextValue 51 5 '  +2.95%  '
